$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 42, shifting existing rows 42-128 down to 43-129.
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new weekly record.
$r = 42
$ws.Cells.Item($r, 1).Value  = 8
$ws.Cells.Item($r, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item($r, 3).Value  = "Coquimbo"
$ws.Cells.Item($r, 4).Value  = 44725
$ws.Cells.Item($r, 5).Value  = 4
$ws.Cells.Item($r, 6).Value  = "Fruta"
$ws.Cells.Item($r, 7).Value  = 100109
$ws.Cells.Item($r, 8).Value  = "Uva"
$ws.Cells.Item($r, 9).Value  = 100109001
$ws.Cells.Item($r, 10).Value = "Uva"
$ws.Cells.Item($r, 11).Value = "Red Globe"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 9500
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 9750
$ws.Cells.Item($r, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($r, 19).Value = 542
$ws.Cells.Item($r, 20).Value = 18
